# "Generate Report for Archive"
#
# b3943904-00b7-48d1-bc7d-ef98103384fe.md and
# bb2ca2c4-ddbe-436c-8253-9617ecc85977.md moved from "Ready for handoff" to
# "In Translation", so the generated report re-sorted the rows: the two
# in-progress files now sit above 87b58ab9-713d-404b-96d0-a38a3bbb52d0.md
# (which stays "Ready for handoff" and drops to the bottom).
#
# Helper: re-point every hyperlink already on a worksheet to a (possibly
# new) display string while preserving its original target URL and anchor
# cell. Hyperlinks.Delete() only works collection-wide in this host, so we
# snapshot {ref -> url} first, wipe them, write the new cell values, and
# recreate each hyperlink against its original URL.
function Set-RowValues($ws, $values) {
    foreach ($ref in $values.Keys) {
        $ws.Range($ref).Value = $values[$ref]
    }
}

function Rebuild-Hyperlinks($ws, $displays) {
    $urls = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $ref = $hl.Range.Address($false, $false)
        $urls[$ref] = $hl.Address
    }
    $ws.Hyperlinks.Delete()
    foreach ($ref in $urls.Keys) {
        $disp = $displays[$ref]
        $ws.Hyperlinks.Add($ws.Range($ref), $urls[$ref], "", "", $disp)
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------- Overview
$ws1 = $wb.Worksheets.Item("Overview")

Set-RowValues $ws1 @{
    "B3" = "In Translation"
    "C3" = "In Translation"
    "D3" = "2016-25-11 14:25:49"

    "B4" = "In Translation"
    "C4" = "In Translation"
    "D4" = "2016-25-11 14:25:49"

    "B5" = "Ready for handoff"
    "C5" = "Ready for handoff"
    "D5" = "2016-24-11 14:24:45"
}

Rebuild-Hyperlinks $ws1 @{
    "A2" = "8135d6b4-f305-45f0-b04b-dcbba50a3d30.md"
    "A3" = "b3943904-00b7-48d1-bc7d-ef98103384fe.md"
    "A4" = "bb2ca2c4-ddbe-436c-8253-9617ecc85977.md"
    "A5" = "87b58ab9-713d-404b-96d0-a38a3bbb52d0.md"
}

# ------------------------------------------------------------------ zh-cn
$ws2 = $wb.Worksheets.Item("zh-cn")

Set-RowValues $ws2 @{
    "C3" = "In Translation"
    "E3" = "2016-03-11 14:25:46"
    "H3" = "0001-01-01 00:00:00"

    "C4" = "In Translation"
    "E4" = "2016-03-11 14:25:46"
    "H4" = "0001-01-01 00:00:00"

    "C5" = "Ready for handoff"
    "E5" = "2016-03-11 14:24:42"
    "H5" = "0001-01-01 00:00:00"
}

Rebuild-Hyperlinks $ws2 @{
    "A2" = "8135d6b4-f305-45f0-b04b-dcbba50a3d30.md"
    "B2" = ".md"
    "D2" = "8135d6b4-f305-45f0-b04b-dcbba50a3d30.08b003844bf439e9423fabc185d486aedb37af0b.zh-cn.xlf"
    "F2" = "8135d6b4-f305-45f0-b04b-dcbba50a3d30.md"
    "G2" = "8135d6b4-f305-45f0-b04b-dcbba50a3d30.08b003844bf439e9423fabc185d486aedb37af0b.zh-cn.xlf"

    "A3" = "b3943904-00b7-48d1-bc7d-ef98103384fe.md"
    "B3" = ".md"
    "D3" = "b3943904-00b7-48d1-bc7d-ef98103384fe.50bcf373b99946182228286ac5c6031ea0bbc7a5.zh-cn.xlf"

    "A4" = "bb2ca2c4-ddbe-436c-8253-9617ecc85977.md"
    "B4" = ".md"
    "D4" = "bb2ca2c4-ddbe-436c-8253-9617ecc85977.e4dd121e6bb0cb618e1951a850a6b3a0fbf011fe.zh-cn.xlf"

    "A5" = "87b58ab9-713d-404b-96d0-a38a3bbb52d0.md"
    "B5" = ".md"
    "D5" = "87b58ab9-713d-404b-96d0-a38a3bbb52d0.2993d00fe1daacb52c128f438d89ecd4cd7f0a7a.zh-cn.xlf"
}

# ------------------------------------------------------------------ de-de
$ws3 = $wb.Worksheets.Item("de-de")

Set-RowValues $ws3 @{
    "C3" = "In Translation"
    "E3" = "2016-03-11 14:25:49"
    "H3" = "0001-01-01 00:00:00"

    "C4" = "In Translation"
    "E4" = "2016-03-11 14:25:49"
    "H4" = "0001-01-01 00:00:00"

    "C5" = "Ready for handoff"
    "E5" = "2016-03-11 14:24:45"
    "H5" = "0001-01-01 00:00:00"
}

Rebuild-Hyperlinks $ws3 @{
    "A2" = "8135d6b4-f305-45f0-b04b-dcbba50a3d30.md"
    "B2" = ".md"
    "D2" = "8135d6b4-f305-45f0-b04b-dcbba50a3d30.08b003844bf439e9423fabc185d486aedb37af0b.de-de.xlf"
    "F2" = "8135d6b4-f305-45f0-b04b-dcbba50a3d30.md"
    "G2" = "8135d6b4-f305-45f0-b04b-dcbba50a3d30.08b003844bf439e9423fabc185d486aedb37af0b.de-de.xlf"

    "A3" = "b3943904-00b7-48d1-bc7d-ef98103384fe.md"
    "B3" = ".md"
    "D3" = "b3943904-00b7-48d1-bc7d-ef98103384fe.50bcf373b99946182228286ac5c6031ea0bbc7a5.de-de.xlf"

    "A4" = "bb2ca2c4-ddbe-436c-8253-9617ecc85977.md"
    "B4" = ".md"
    "D4" = "bb2ca2c4-ddbe-436c-8253-9617ecc85977.e4dd121e6bb0cb618e1951a850a6b3a0fbf011fe.de-de.xlf"

    "A5" = "87b58ab9-713d-404b-96d0-a38a3bbb52d0.md"
    "B5" = ".md"
    "D5" = "87b58ab9-713d-404b-96d0-a38a3bbb52d0.2993d00fe1daacb52c128f438d89ecd4cd7f0a7a.de-de.xlf"
}
